# Add a new date column (CI, 2024/12/04) to the right of the existing
# last date column (CH, 2024/12/03) on the sheet's single data table.
#
# Row 1 holds the date header (text). Rows 2-53 hold a numeric "back ratio"
# per machine, highlighted with one of three pre-existing cell styles
# depending on the value:
#   - style "2" (yellow fill)     when value <  125.0
#   - style "3" (light-blue fill) when 125.0 <= value < 140.0
#   - style "1" (no fill)         when value >= 140.0
# (derived from the existing sheet: every s="2" cell is in [92.0, 124.9],
# every s="3" cell is in [125.0, 139.9], every s="1" numeric cell is >= 140).
#
# Rather than poke Font/Interior properties directly (which bakes in new
# style entries because the COM layer materialises extra attributes such
# as an explicit theme color), each new cell is produced by copying an
# existing same-style cell (value + format) onto it and then overwriting
# just the value — so the new cells land on the exact same style indexes
# already used elsewhere in the sheet instead of creating new ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcCol = 86   # CH (2024/12/03)
$dstCol = 87   # CI (2024/12/04)

# New column gets the same width as its neighbour.
$ws.Columns.Item($dstCol).ColumnWidth = $ws.Columns.Item($srcCol).ColumnWidth

# --- Row 1 (header, text date) -------------------------------------------
# Copying CH1 first (value + format) puts CI1 on CH1's exact style ("1").
# Forcing text format before assigning the literal keeps Excel from
# parsing "2024/12/04" as a date serial; the immediately-following
# "paste formats" from CH1 then restores the plain style ("1") that the
# text-format tweak had bumped to a throwaway one-off style, without
# disturbing the text value just assigned.
$headerDst = $ws.Cells.Item(1, $dstCol)
$headerSrc = $ws.Cells.Item(1, $srcCol)
$headerDst.NumberFormat = "@"
$headerDst.Value = "2024/12/04"
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rows 2-53 (numeric values) -------------------------------------------
# Reference cells already on each of the three styles, found in column CH
# itself, used as copy sources so no new style entries get created.
$refNormal = $ws.Cells.Item(3, $srcCol)   # CH3  -> style "1" (no fill)
$refYellow = $ws.Cells.Item(2, $srcCol)   # CH2  -> style "2" (yellow fill)
$refBlue   = $ws.Cells.Item(21, $srcCol)  # CH21 -> style "3" (light-blue fill)

$newValues = @{
    2  = 136.3
    3  = 118.5
    4  = 139.8
    5  = 121.7
    6  = 199.4
    7  = 163.4
    8  = 155.1
    9  = 188.5
    10 = 131.7
    11 = 156.9
    12 = 150.8
    13 = 147.4
    14 = 117.5
    15 = 191.3
    16 = 173.8
    17 = 167.5
    18 = 146.6
    19 = 133.4
    20 = 160.3
    21 = 137.3
    22 = 241.3
    23 = 140
    24 = 136.7
    25 = 174.6
    26 = 218.5
    27 = 135.8
    28 = 180.2
    29 = 227.1
    30 = 337
    31 = 166.7
    32 = 129.8
    33 = 132.9
    34 = 246.5
    35 = 124.8
    36 = 143.2
    37 = 428.3
    38 = 173.1
    39 = 157.9
    40 = 192.6
    41 = 138.4
    42 = 134.8
    43 = 157.7
    44 = 136.9
    45 = 258.8
    46 = 147.5
    47 = 128.6
    48 = 201.8
    49 = 205
    50 = 240
    51 = 196.4
    52 = 146.2
    53 = 164.8
}

for ($row = 2; $row -le 53; $row++) {
    $value = $newValues[$row]

    if ($value -lt 125.0) {
        $ref = $refYellow
    } elseif ($value -lt 140.0) {
        $ref = $refBlue
    } else {
        $ref = $refNormal
    }

    $dst = $ws.Cells.Item($row, $dstCol)
    $ref.Copy($dst)
    $dst.Value = $value
}

$excel.CutCopyMode = 0
